$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Move the Gantt chart's visible week back to week 1 (Display_Week, named range -> D4).
# All the date header formulas (G4/G5 .. JM4/JS5 etc.) recompute automatically
# off this single input cell.
$ws.Range("D4").Value = 1

# Bump the "In Progress" sub-task completion percentage; "Complete" (C22) is a
# formula average of C23:C24 and recalculates on its own.
$ws.Range("C23").Value = 0.7

# View: zoom back out and move the active selection to D5.
$excel.ActiveWindow.Zoom = 85
$ws.Range("D5").Select()

# Print setup: drop the explicit paper size, shrink the scale way down, and
# switch from landscape to portrait, while keeping the "no row limit" setting.
$ws.PageSetup.PaperSize = 0
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Zoom = 10
$ws.PageSetup.FitToPagesTall = $False
